# Daily attendance processing - 2025-10-22 10:50:05
# Reorders the "Recorded By" (column G) audit-trail names for rows where the
# contributor list changed ordering (system/System/backup@backdoor.com entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value2 = "system, backup@backdoor.com, System"
    } elseif ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    } elseif ($val -eq "System, backup@backdoor.com") {
        $cell.Value2 = "backup@backdoor.com, System"
    }
}
